$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Sheet1" (xl/worksheets/sheet2.xml) - add a new row with a
# single space value, used as a "unwanted space" field.
# -----------------------------------------------------------------
$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Range("A6").Value = " "
$wsSheet1.Range("A5").Select()

# -----------------------------------------------------------------
# Sheet "Sheet2" (xl/worksheets/sheet3.xml) - extend the table with
# extra columns (K..S) for percentage / hub-collection data, and add
# the new "hub collection" share + button label columns (R, S).
# -----------------------------------------------------------------
$wsSheet2 = $wb.Worksheets.Item("Sheet2")

# New column widths for the added columns.
$wsSheet2.Columns.Item(8).ColumnWidth = 17
$wsSheet2.Columns.Item(18).ColumnWidth = 17.083333333333336

# Row 1 header cells.
$wsSheet2.Range("J1").Value = "DELAY"
$wsSheet2.Range("K1").Value = "SWIPE"
$wsSheet2.Range("L1").Value = "SWIPE"
$wsSheet2.Range("M1").Value = "SWIPE"
$wsSheet2.Range("N1").Value = "SWIPE"
$wsSheet2.Range("O1").Value = "SWIPE"
$wsSheet2.Range("P1").Value = "SWIPE"
$wsSheet2.Range("Q1").Value = "SWIPE"
$wsSheet2.Range("R1").Value = "CLICK"
$wsSheet2.Range("S1").Value = "CLICK"

# Row 2 header / sub-header cells.
$wsSheet2.Range("I2").ClearContents()
$wsSheet2.Range("J2").Value = 5
$wsSheet2.Range("K2").Value = "RIGHT"
$wsSheet2.Range("L2").Value = "RIGHT"
$wsSheet2.Range("M2").Value = "RIGHT"
$wsSheet2.Range("N2").Value = "RIGHT"
$wsSheet2.Range("O2").Value = "RIGHT"
$wsSheet2.Range("P2").Value = "RIGHT"
$wsSheet2.Range("Q2").Value = "RIGHT"

# New "hub collection" columns use a dedicated PingFang SC font. We
# build the formatting through a temporary named style so only a
# single new font / cell format gets registered, then drop the named
# style again (keeping the direct formatting on the cell).
$hubStyle = $wb.Styles.Add("HubCollectionStyle")
$hubStyle.Font.Name = "PingFang SC"
$wsSheet2.Range("R2").Style = "HubCollectionStyle"
$wb.Styles.Item("HubCollectionStyle").Delete()

$wsSheet2.Range("R2").Value = "ส่วนแบ่งhubCollection"
$wsSheet2.Range("S2").Value = "เพิ่มเปอร์เซ็นต์hubcollection_button"

# Row 3 - the "Item_0" label moves from I2 down to I3.
$wsSheet2.Range("I3").Value = "Item_0"

# Row height for row 2 grew slightly once the extra content was added.
$wsSheet2.Rows.Item(2).RowHeight = 19

$wsSheet2.Range("T1").Select()
$wsSheet2.Activate()
